$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$th = $nm.Theme
Write-Host "Theme.Name=$($th.Name)"
try { Write-Host "ThemeColorScheme=$($th.ThemeColorScheme)" } catch { Write-Host "err: $_" }
try { Write-Host "Save=$($th.Save)" } catch { Write-Host "err save: $_" }

$sm = $p.SlideMaster
$th2 = $sm.Theme
Write-Host "SlideMaster Theme.Name=$($th2.Name)"
